# Auto-generated: apply numeric updates to Pandaemonium_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Cells.Item(69, 8).Value = 4818
$ws.Cells.Item(69, 9).Value = 4731.5
$ws.Cells.Item(69, 10).Value = 4933.3335
$ws.Cells.Item(69, 11).Value = 14194.5
$ws.Cells.Item(69, 12).Value = 14800.0005
$ws.Cells.Item(69, 13).Value = -13320.5
$ws.Cells.Item(69, 14).Value = -16548.0005

# Row 72
$ws.Cells.Item(72, 8).Value = 4818
$ws.Cells.Item(72, 9).Value = 4731.5
$ws.Cells.Item(72, 10).Value = 4933.3335
$ws.Cells.Item(72, 11).Value = 42583.5
$ws.Cells.Item(72, 12).Value = 44400.0015
$ws.Cells.Item(72, 13).Value = -38215.5
$ws.Cells.Item(72, 14).Value = -53136.0015

# Row 112
$ws.Cells.Item(112, 8).Value = 1050.091
$ws.Cells.Item(112, 10).Value = 1076.762
$ws.Cells.Item(112, 12).Value = 3230.286
$ws.Cells.Item(112, 14).Value = -5446.286

# Row 138
$ws.Cells.Item(138, 8).Value = 3777.1372
$ws.Cells.Item(138, 9).Value = 1509.4117
$ws.Cells.Item(138, 10).Value = 4911
$ws.Cells.Item(138, 11).Value = 4528.2351
$ws.Cells.Item(138, 12).Value = 14733
$ws.Cells.Item(138, 13).Value = 611.7649000000001
$ws.Cells.Item(138, 14).Value = -25013

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Cells.Item(63, 8).Value = 3932.5
$ws.Cells.Item(63, 9).Value = 3721
$ws.Cells.Item(63, 11).Value = 3721
$ws.Cells.Item(63, 13).Value = -3035

# Row 66
$ws.Cells.Item(66, 8).Value = 3932.5
$ws.Cells.Item(66, 9).Value = 3721
$ws.Cells.Item(66, 11).Value = 18605
$ws.Cells.Item(66, 13).Value = -15173

# Row 74
$ws.Cells.Item(74, 8).Value = 5682.865
$ws.Cells.Item(74, 9).Value = 4093.628
$ws.Cells.Item(74, 10).Value = 13275.889
$ws.Cells.Item(74, 11).Value = 4093.628
$ws.Cells.Item(74, 12).Value = 13275.889
$ws.Cells.Item(74, 13).Value = -3219.628
$ws.Cells.Item(74, 14).Value = -15023.889

# Row 77
$ws.Cells.Item(77, 8).Value = 5682.865
$ws.Cells.Item(77, 9).Value = 4093.628
$ws.Cells.Item(77, 10).Value = 13275.889
$ws.Cells.Item(77, 11).Value = 20468.14
$ws.Cells.Item(77, 12).Value = 66379.44499999999
$ws.Cells.Item(77, 13).Value = -16100.14
$ws.Cells.Item(77, 14).Value = -75115.44499999999

# Row 102
$ws.Cells.Item(102, 8).Value = 1483364.9
$ws.Cells.Item(102, 9).Value = 1765072.5
$ws.Cells.Item(102, 11).Value = 1765072.5
$ws.Cells.Item(102, 13).Value = -1763450.5

# Row 122
$ws.Cells.Item(122, 8).Value = 2976868.8
$ws.Cells.Item(122, 9).Value = 628
$ws.Cells.Item(122, 10).Value = 41668000
$ws.Cells.Item(122, 11).Value = 1884
$ws.Cells.Item(122, 12).Value = 125004000
$ws.Cells.Item(122, 13).Value = 566
$ws.Cells.Item(122, 14).Value = -125008900

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 1898.9166
$ws.Cells.Item(94, 9).Value = 1737.125
$ws.Cells.Item(94, 10).Value = 2222.5
$ws.Cells.Item(94, 11).Value = 1737.125
$ws.Cells.Item(94, 12).Value = 2222.5
$ws.Cells.Item(94, 13).Value = -1286.125
$ws.Cells.Item(94, 14).Value = -3124.5

# Row 99
$ws.Cells.Item(99, 8).Value = 1326.8928
$ws.Cells.Item(99, 9).Value = 1168.5555
$ws.Cells.Item(99, 10).Value = 1611.9
$ws.Cells.Item(99, 11).Value = 1168.5555
$ws.Cells.Item(99, 12).Value = 1611.9
$ws.Cells.Item(99, 13).Value = 329.4445000000001
$ws.Cells.Item(99, 14).Value = -4607.9

# Row 105
$ws.Cells.Item(105, 8).Value = 764467.4
$ws.Cells.Item(105, 9).Value = 803401.3
$ws.Cells.Item(105, 10).Value = 5255.5
$ws.Cells.Item(105, 11).Value = 803401.3
$ws.Cells.Item(105, 12).Value = 5255.5
$ws.Cells.Item(105, 13).Value = -801654.3
$ws.Cells.Item(105, 14).Value = -8749.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2965.959
$ws.Cells.Item(31, 9).Value = 2727.7212
$ws.Cells.Item(31, 10).Value = 4177
$ws.Cells.Item(31, 11).Value = 2727.7212
$ws.Cells.Item(31, 12).Value = 4177
$ws.Cells.Item(31, 13).Value = -2432.7212
$ws.Cells.Item(31, 14).Value = -4767

# Row 34
$ws.Cells.Item(34, 8).Value = 2965.959
$ws.Cells.Item(34, 9).Value = 2727.7212
$ws.Cells.Item(34, 10).Value = 4177
$ws.Cells.Item(34, 11).Value = 2727.7212
$ws.Cells.Item(34, 12).Value = 4177
$ws.Cells.Item(34, 13).Value = -2525.7212
$ws.Cells.Item(34, 14).Value = -4581

# Row 58
$ws.Cells.Item(58, 8).Value = 1573.527
$ws.Cells.Item(58, 9).Value = 915.92725
$ws.Cells.Item(58, 10).Value = 3477.1052
$ws.Cells.Item(58, 11).Value = 915.92725
$ws.Cells.Item(58, 12).Value = 3477.1052
$ws.Cells.Item(58, 13).Value = -712.92725
$ws.Cells.Item(58, 14).Value = -3883.1052

# Row 99
$ws.Cells.Item(99, 8).Value = 3528.1667
$ws.Cells.Item(99, 9).Value = 1600
$ws.Cells.Item(99, 10).Value = 3913.8
$ws.Cells.Item(99, 11).Value = 1600
$ws.Cells.Item(99, 12).Value = 3913.8
$ws.Cells.Item(99, 13).Value = -102
$ws.Cells.Item(99, 14).Value = -6909.8

# Row 126
$ws.Cells.Item(126, 8).Value = 3528.1667
$ws.Cells.Item(126, 9).Value = 1600
$ws.Cells.Item(126, 10).Value = 3913.8
$ws.Cells.Item(126, 11).Value = 4800
$ws.Cells.Item(126, 12).Value = 11741.4
$ws.Cells.Item(126, 13).Value = -2330
$ws.Cells.Item(126, 14).Value = -16681.4

# Row 132
$ws.Cells.Item(132, 8).Value = 6626.5
$ws.Cells.Item(132, 9).Value = 11198.25
$ws.Cells.Item(132, 10).Value = 3578.6667
$ws.Cells.Item(132, 11).Value = 33594.75
$ws.Cells.Item(132, 12).Value = 10736.0001
$ws.Cells.Item(132, 13).Value = -31064.75
$ws.Cells.Item(132, 14).Value = -15796.0001

# Row 134
$ws.Cells.Item(134, 8).Value = 2874.074
$ws.Cells.Item(134, 9).Value = 2056.1738
$ws.Cells.Item(134, 10).Value = 3480.9033
$ws.Cells.Item(134, 11).Value = 6168.5214
$ws.Cells.Item(134, 12).Value = 10442.7099
$ws.Cells.Item(134, 13).Value = -3633.5214
$ws.Cells.Item(134, 14).Value = -15512.7099

# Row 136
$ws.Cells.Item(136, 8).Value = 1573.527
$ws.Cells.Item(136, 9).Value = 915.92725
$ws.Cells.Item(136, 10).Value = 3477.1052
$ws.Cells.Item(136, 11).Value = 2747.78175
$ws.Cells.Item(136, 12).Value = 10431.3156
$ws.Cells.Item(136, 13).Value = -197.7817500000001
$ws.Cells.Item(136, 14).Value = -15531.3156

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Cells.Item(92, 8).Value = 686.3333
$ws.Cells.Item(92, 10).Value = 742
$ws.Cells.Item(92, 12).Value = 2226
$ws.Cells.Item(92, 14).Value = -4722

# Row 107
$ws.Cells.Item(107, 8).Value = 1659.6
$ws.Cells.Item(107, 9).Value = 325
$ws.Cells.Item(107, 10).Value = 2549.3333
$ws.Cells.Item(107, 11).Value = 975
$ws.Cells.Item(107, 12).Value = 7647.999899999999
$ws.Cells.Item(107, 13).Value = 945
$ws.Cells.Item(107, 14).Value = -11487.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3722.5757
$ws.Cells.Item(40, 9).Value = 3494.8333
$ws.Cells.Item(40, 10).Value = 6000
$ws.Cells.Item(40, 11).Value = 3494.8333
$ws.Cells.Item(40, 12).Value = 6000
$ws.Cells.Item(40, 13).Value = -3358.8333
$ws.Cells.Item(40, 14).Value = -6272

# Row 88
$ws.Cells.Item(88, 8).Value = 39911.832
$ws.Cells.Item(88, 9).Value = 40085.5
$ws.Cells.Item(88, 10).Value = 39825
$ws.Cells.Item(88, 11).Value = 40085.5
$ws.Cells.Item(88, 12).Value = 39825
$ws.Cells.Item(88, 13).Value = -39657.5
$ws.Cells.Item(88, 14).Value = -40681

# Row 91
$ws.Cells.Item(91, 8).Value = 39911.832
$ws.Cells.Item(91, 9).Value = 40085.5
$ws.Cells.Item(91, 10).Value = 39825
$ws.Cells.Item(91, 11).Value = 40085.5
$ws.Cells.Item(91, 12).Value = 39825
$ws.Cells.Item(91, 13).Value = -38603.5
$ws.Cells.Item(91, 14).Value = -42789

# Row 101
$ws.Cells.Item(101, 8).Value = 5833
$ws.Cells.Item(101, 10).Value = 5833
$ws.Cells.Item(101, 12).Value = 5833
$ws.Cells.Item(101, 14).Value = -12323

# Row 132
$ws.Cells.Item(132, 8).Value = 4055.2546
$ws.Cells.Item(132, 9).Value = 4064.9268
$ws.Cells.Item(132, 10).Value = 4026.9285
$ws.Cells.Item(132, 11).Value = 12194.7804
$ws.Cells.Item(132, 12).Value = 12080.7855
$ws.Cells.Item(132, 13).Value = -9664.7804
$ws.Cells.Item(132, 14).Value = -17140.7855

# Row 136
$ws.Cells.Item(136, 8).Value = 4465.9214
$ws.Cells.Item(136, 9).Value = 2599.9143
$ws.Cells.Item(136, 11).Value = 7799.742899999999
$ws.Cells.Item(136, 13).Value = -5249.742899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Cells.Item(103, 8).Value = 60602
$ws.Cells.Item(103, 12).Value = 60602
$ws.Cells.Item(103, 14).Value = -62946
